$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the existing hyperlink before restructuring rows so that the
#     underlying cell styling can be cleanly reassigned afterwards. ---
$ws.Range("D8").Hyperlinks.Delete()

# --- Insert a new row for author #2's "id" field (650798), mirroring the
#     "id" row that already exists for author #1 (row 3). This pushes the
#     remaining "authors" detail rows (name/first_name/last_name/email) for
#     author #3 down by one row. The new row automatically inherits the
#     formatting of the row above it (row 4), which already matches the
#     styling used throughout the "authors" block. ---
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "authors"
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = "id"
$ws.Cells.Item(5, 4).Value = 650798

# --- Remove the now-redundant "orcid_id" detail row, which sat at row 9
#     before the insert above and has shifted down to row 10. ---
$ws.Rows.Item(10).Delete()

# --- Re-create the hyperlink on the email cell, which is now row 9. ---
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:thebard@hotmail.com") | Out-Null

# Restore the plain (non-hyperlink) style on the old hyperlink cell (D8,
# now holding the "last_name" value) and make sure D9 uses the Hyperlink
# style (re-applying it re-uses the workbook's existing "Hyperlink" cell
# style rather than leaving the ad-hoc one that Hyperlinks.Add() applies).
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D9").Style = "Hyperlink"

# Update the active selection to match the edited workbook.
$ws.Range("H23").Select()
